$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-HeaderStyle($r) {
    $r.Font.Bold = $true
    $r.Borders.LineStyle = 1
    $r.HorizontalAlignment = -4108
    $r.VerticalAlignment = -4160
}

# New header cells (row 1)
$ws.Range("AH1").Value = 'isRound'
Set-HeaderStyle $ws.Range("AH1")
$ws.Range("AI1").Value = 'highRisk'
Set-HeaderStyle $ws.Range("AI1")
$ws.Range("AJ1").Value = 'ComplFATF'
Set-HeaderStyle $ws.Range("AJ1")

# Row 2
$ws.Range("A2").Value = 'example_1.txt'
Set-HeaderStyle $ws.Range("A2")
$ws.Range("B2").Value = 44641
$ws.Range("B2").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C2").Value = 'MT103 0001'
$ws.Range("D2").Value = '/INS/THIS IS A PAYMENT FOR TUNA SUPPLY '''
$ws.Range("E2").Value = 'USD'
$cell = $ws.Range("F2")
$cell.NumberFormat = "@"
$cell.Value = '5000'
$cell.ClearFormats()
$ws.Range("G2").Value = 'CRED'
$ws.Range("K2").Value = 'COMMERZBANK AG'
$ws.Range("L2").Value = 'COMMERZBANK'
$ws.Range("N2").Value = 'AG'
$ws.Range("O2").Value = 'HAMBURG, GERMANY'
$ws.Range("P2").Value = 'DE'
$ws.Range("Q2").Value = 'DE98765432101234567890'
$ws.Range("R2").Value = 'XXX'
$ws.Range("S2").Value = 'COBADEHHXXX'
$ws.Range("T2").Value = 'METRO BANK PLC'
$ws.Range("U2").Value = 'LONDON, UNITED KINGDOM'
$ws.Range("V2").Value = 'Not found'
$ws.Range("W2").Value = 'NORDFISCH GMBH'
$ws.Range("X2").Value = 'NORDFISCH'
$ws.Range("Z2").Value = 'GMBH'
$ws.Range("AA2").Value = 'BODENSEE STR. 226 22761 HAMBURG GERMANY'
$ws.Range("AB2").Value = 'GB'
$ws.Range("AC2").Value = 'GB57METR12345678901234'
$ws.Range("AD2").Value = 'XXX'
$ws.Range("AE2").Value = 'HBUKGB4BXXX'
$ws.Range("AF2").Value = 'Not found'
$ws.Range("AG2").Value = 'Not found'
$ws.Range("AH2").Value = $true
$ws.Range("AI2").Value = $false
$ws.Range("AJ2").Value = $false

# Row 3
$ws.Range("A3").Value = 'example_2.txt'
Set-HeaderStyle $ws.Range("A3")
$ws.Range("B3").Value = 44641
$ws.Range("B3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C3").Value = 'MT103 0001'
$ws.Range("D3").Value = '/MSG/PAYMENT FOR GOODS "'
$ws.Range("E3").Value = 'USD'
$cell = $ws.Range("F3")
$cell.NumberFormat = "@"
$cell.Value = '10000'
$cell.ClearFormats()
$ws.Range("G3").Value = 'CRED'
$ws.Range("K3").Value = 'ABC INDUSTRIES'
$ws.Range("L3").Value = 'ABC'
$ws.Range("N3").Value = 'INDUSTRIES'
$ws.Range("O3").Value = '123 MAIN STREET NEW YORK, NY 10001 UNITED STATES'
$ws.Range("P3").Value = 'US'
$ws.Range("Q3").Value = 'US12345678901234567890'
$ws.Range("R3").Value = 'XXX'
$ws.Range("S3").Value = 'ABCBUS33XXX'
$ws.Range("T3").Value = 'HSBC HONG KONG'
$ws.Range("U3").Value = 'HONG KONG'
$ws.Range("V3").Value = 'ICBKCNBJGZU'
$ws.Range("W3").Value = 'XYZ SUPPLIERS'
$ws.Range("X3").Value = 'XYZ'
$ws.Range("Z3").Value = 'SUPPLIERS'
$ws.Range("AA3").Value = '123 HUANGPU ROAD SHANGHAI, CHINA'
$ws.Range("AB3").Value = 'CN'
$ws.Range("AC3").Value = 'CN123456789012345678'
$ws.Range("AD3").Value = 'Not found'
$ws.Range("AE3").Value = 'CITIUS33'
$ws.Range("AF3").Value = 'CITIBANK HONG KONG'
$ws.Range("AG3").Value = 'CENTRAL, HONG KONG'
$ws.Range("AH3").Value = $true
$ws.Range("AI3").Value = $false
$ws.Range("AJ3").Value = $false

# Row 4
$ws.Range("A4").Value = 'example_3.txt'
Set-HeaderStyle $ws.Range("A4")
$ws.Range("B4").Value = 44641
$ws.Range("B4").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("C4").Value = 'MT103 0001'
$ws.Range("D4").Value = 'INV NO. 12345 REF. 98765 SUPPLY OF GOODS AS PER PURCHASE ORDER NO. 54321 '''
$ws.Range("E4").Value = 'USD'
$cell = $ws.Range("F4")
$cell.NumberFormat = "@"
$cell.Value = '9899'
$cell.ClearFormats()
$ws.Range("G4").Value = 'CRED'
$ws.Range("K4").Value = 'ABC SUPPLIERS BV'
$ws.Range("L4").Value = 'ABC'
$ws.Range("M4").Value = 'SUPPLIERS'
$ws.Range("N4").Value = 'BV'
$ws.Range("O4").Value = 'AMSTERDAM, NETHERLANDS'
$ws.Range("P4").Value = 'NL'
$ws.Range("Q4").Value = 'NL20ABNA0404875234'
$ws.Range("R4").Value = 'Not found'
$ws.Range("S4").Value = 'ABNANL2A'
$ws.Range("T4").Value = 'Not found'
$ws.Range("U4").Value = 'Not found'
$ws.Range("V4").Value = 'SCBLGB2LXXX'
$ws.Range("W4").Value = 'AFRICAN EXPORT-IMPORT BANK'
$ws.Range("X4").Value = 'AFRICAN'
$ws.Range("Y4").Value = 'EXPORT-IMPORT'
$ws.Range("Z4").Value = 'BANK'
$ws.Range("AA4").Value = 'LAGOS, NIGERIA XYZ ENTERPRISES LTD LAGOS, NIGERIA'
$ws.Range("AB4").Value = 'Not found'
$ws.Range("AC4").Value = 'Not found'
$ws.Range("AD4").Value = 'XXX'
$ws.Range("AE4").Value = 'PASSNGLAXXX'
$ws.Range("AF4").Value = 'Not found'
$ws.Range("AG4").Value = 'Not found'
$ws.Range("AH4").Value = $false
$ws.Range("AI4").Value = $false
$ws.Range("AJ4").Value = $false
